# Trial_Sentencing_Template.docx -- "Updated trial sentencing template to
# remove auto victim statement."
#
# The paragraph that lists who the Court heard statements from currently
# always mentions the Victim Advocate. This script:
#   1. Splits the run so the leading "{{'\n\n'}}" marker stays on its own,
#      unchanged run.
#   2. Rewrites the sentence so the victim-advocate clause is now gated by
#      {% if victim_statements is true %} ... {% endif %}, and both the
#      victim-advocate and defense-counsel clauses lead with a comma
#      (instead of the comma trailing the clause).
#   3. Adds a trailing space to the lone "}" run that closes the
#      surrounding {% if defense_counsel_waived ... %...} template tag.

$d = $word.ActiveDocument

$oldSentence = "{{‘\n\n’}}Prior to sentencing, the Court heard statements from the Prosecutor, Victim Advocate on behalf of the victim, {% if defense_counsel_waived is false %}Defense Counsel {% endif %"

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found) {
    throw "Could not locate the target sentence in word/document.xml"
}

$matchRange = $find.Parent
$matchStart = $matchRange.Start
$matchEnd = $matchRange.End

# "{{'\n\n'}}" (the template newline marker) is 10 characters long and
# must be left alone, in its own run.
$markerLen = 10
$splitPoint = $matchStart + $markerLen

$newBody = "Prior to sentencing, the Court heard statements from the Prosecutor{% if victim_statements is true %}, Victim Advocate on behalf of the victim{% endif %}{% if defense_counsel_waived is false %}, Defense Counsel{% endif %"

# Rewrite the sentence body (everything after the marker, up to but not
# including the closing "}" run).
$bodyRange = $d.Range($splitPoint, $matchEnd)
$bodyRange.Text = $newBody

$newMatchEnd = $splitPoint + $newBody.Length

# Force the rewritten body to sit in its own run, distinct from the
# "{{'\n\n'}}" marker run, by toggling a character property on and back
# off -- this splits the run while leaving its final formatting (and
# thus the neighbouring run's formatting) unchanged.
$bodyRange2 = $d.Range($splitPoint, $newMatchEnd)
$bodyRange2.Font.Bold = $true
$bodyRange2.Font.Bold = $false

# The very next character is the lone "}" that closes the
# {% if defense_counsel_waived ... %...} tag; give it a trailing space.
$closeRange = $d.Range($newMatchEnd, $newMatchEnd + 1)
$closeRange.Text = "} "
